$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price column (D) keeps its original plain-text storage,
# since several prices (e.g. "1.003") would otherwise be parsed as numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.776.66'
$ws.Range("E2").Value = '  -0.88%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.627.45'
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.35'
$ws.Range("E5").Value = '  +0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5056'
$ws.Range("E6").Value = '  -0.38%  '
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2571'
$ws.Range("E8").Value = '  -0.39%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06404'
$ws.Range("E9").Value = '  +0.69%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.39'
$ws.Range("E10").Value = '  -2.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07799'
$ws.Range("E11").Value = '  +0.82%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.255'
$ws.Range("E12").Value = '  -1.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.629.06'
$ws.Range("E13").Value = '  -0.50%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.851.68'
$ws.Range("E14").Value = '  -0.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5582'
$ws.Range("E15").Value = '  +1.91%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '62.95'
$ws.Range("E16").Value = '  -2.29%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0₅7540'
$ws.Range("E17").Value = '  -2.79%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.791.60'
$ws.Range("E18").Value = '  -0.91%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  +0.25%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '193.89'
$ws.Range("E20").Value = '  -1.71%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.317'
$ws.Range("E21").Value = '  -3.24%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.811'
$ws.Range("E22").Value = '  -1.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.987'
$ws.Range("E23").Value = '  -2.47%  '
$ws.Range("E24").Value = '  +0.13%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.793'
$ws.Range("E25").Value = '  -5.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '140.84'
$ws.Range("E26").Value = '  -1.31%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1265'
$ws.Range("E27").Value = '  +0.57%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.728'
$ws.Range("E28").Value = '  -2.15%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.38'
$ws.Range("E29").Value = '  -1.42%  '
$ws.Range("E30").Value = '  -0.38%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04862'
$ws.Range("E31").Value = '  -0.80%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.278'
$ws.Range("E32").Value = '  -0.28%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.194'
$ws.Range("E33").Value = '  -0.70%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.556'
$ws.Range("E34").Value = '  -0.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.376'
$ws.Range("E35").Value = '  -0.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8937'
$ws.Range("E36").Value = '  -2.75%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.566'
$ws.Range("E37").Value = '  +0.04%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.128.08'
$ws.Range("E38").Value = '  +2.03%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5470'
$ws.Range("E39").Value = '  -1.56%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.01562'
$ws.Range("E40").Value = '  -0.47%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.9986'
$ws.Range("E41").Value = '  -0.25%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.548'
$ws.Range("E42").Value = '  -1.20%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.7965'
$ws.Range("E43").Value = '  -1.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '97.29'
$ws.Range("E44").Value = '  -1.50%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.779.62'
$ws.Range("E45").Value = '  -0.08%  '
$ws.Range("E46").Value = '  -6.00%  '
$ws.Range("E47").Value = '  -2.18%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '55.16'
$ws.Range("E48").Value = '  -0.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.05054'
$ws.Range("E49").Value = '  -2.74%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.639'
$ws.Range("E50").Value = '  +1.12%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.000'
$ws.Range("E51").Value = '  -0.15%  '
